$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 5.815
$ws.Range("B13").Value = 5.556
$ws.Range("B16").Value = 4.871
$ws.Range("B18").Value = 5.255999999999999
$ws.Range("B20").Value = 6.862
$ws.Range("B26").Value = 5.897
$ws.Range("B27").Value = 5.752000000000001
$ws.Range("B29").Value = 5.633
$ws.Range("B35").Value = 8.301
$ws.Range("B36").Value = 8.079000000000001
$ws.Range("B45").Value = 5.516
$ws.Range("B55").Value = 4.813
$ws.Range("B57").Value = 5.233000000000001
$ws.Range("B69").Value = 5.128
$ws.Range("B76").Value = 6.544999999999999
$ws.Range("B78").Value = 8.221
$ws.Range("B82").Value = 5.456999999999999
$ws.Range("B83").Value = 5.88
$ws.Range("B93").Value = 4.896000000000001
$ws.Range("B97").Value = 4.86
